# Updates cryptos list values (price + 1h volume change) to match
# the scraped snapshot described in the commit diff, plus the
# Monero / LidoDAOToken row swap (rows 25-26).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.644.31"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").Value = "1.844.96"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "260.36"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5279"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.95%  "
$ws.Range("E8").Value = "  -3.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06798"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.79"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7841"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07761"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").Value = "1.840.49"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.07"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.016"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9991"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007928"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").Value = "26.669.79"
$ws.Range("E20").Value = "  +0.82%  "
$ws.Range("D21").Value = "2.081.56"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.611"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.982"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.316"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.17%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.63"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.03%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.218"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.683"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.88%  "
$ws.Range("E28").Value = "  +0.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "110.92"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.204"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08727"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.086"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.80%  "
$ws.Range("E33").Value = "  +1.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7301"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.143"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.858"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.098"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.281"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.98%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01734"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4784"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8992"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.74"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.955"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.08%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4172"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.076"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("E48").Value = "  +1.88%  "
$ws.Range("E49").Value = "  -2.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.79"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8918"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.73%  "
